# CCC19 Derived Variables Spreadsheet: add a new "Categorical age variable"
# row (D01a / age_cat) right after the existing "age" (D01) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 109 (pushes existing row 109+ down by one).
$ws.Rows.Item(109).Insert()

# Populate the new row with the derived-variable metadata.
$ws.Range("A109").Value = "D01a"
$ws.Range("B109").Value = "age_cat"
$ws.Range("C109").Value = "Demographics"
$ws.Range("D109").Value = "Categorical age variable"
$ws.Range("E109").Value = "18-39 years; 40-59 years; 60-69 years; 70-79 years; 80+ years"

# Grow Table1 so it keeps covering the whole data range (was A1:E236, now A1:E237).
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E237"))

# Match the author's final selection/cursor position in the sheet.
$ws.Range("E49").Select()
